# Sean Matthew timesheet (2026-01-26 week) — full-month coverage fix:
#  - replace placeholder client names with the actual simulator output
#  - day 1 becomes a PTO day (no billable client)
#  - hours/rate/total recomputed for every day + all subtotal rows
#  - employee id regenerated

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# ---- Weekly Timesheet sheet (A:Date B:Client C:Hours D:Type E:Rate F:Total) ----

# Row 2 - 2026-01-26 -> PTO day
$wsWeekly.Range("B2").Value = "PTO"
$wsWeekly.Range("C2").Value = 6
$wsWeekly.Range("D2").Value = "PTO"
$wsWeekly.Range("E2").Value = 88
$wsWeekly.Range("F2").Value = 528

# Row 3 - 2026-01-27
$wsWeekly.Range("B3").Value = "Vincent"
$wsWeekly.Range("C3").Value = 6.5
$wsWeekly.Range("D3").Value = "Regular"
$wsWeekly.Range("E3").Value = 88
$wsWeekly.Range("F3").Value = 572

# Row 4 - 2026-01-28
$wsWeekly.Range("B4").Value = "Zygmunt"
$wsWeekly.Range("C4").Value = 6
$wsWeekly.Range("D4").Value = "Regular"
$wsWeekly.Range("E4").Value = 88
$wsWeekly.Range("F4").Value = 528

# Row 5 - 2026-01-29
$wsWeekly.Range("B5").Value = "Ricca"
$wsWeekly.Range("C5").Value = 6.5
$wsWeekly.Range("D5").Value = "Regular"
$wsWeekly.Range("E5").Value = 88
$wsWeekly.Range("F5").Value = 572

# Row 6 - 2026-01-30
$wsWeekly.Range("B6").Value = "Varricchio"
$wsWeekly.Range("C6").Value = 7
$wsWeekly.Range("D6").Value = "Regular"
$wsWeekly.Range("E6").Value = 88
$wsWeekly.Range("F6").Value = 616

# SUBTOTAL row
$wsWeekly.Range("C8").Value = 32
$wsWeekly.Range("D8").Value = "Reg: 32 / OT: 0"
$wsWeekly.Range("F8").Value = 2816

# HOURLY SUBTOTAL row
$wsWeekly.Range("F11").Value = 2816

# GRAND TOTAL row
$wsWeekly.Range("F13").Value = 2816

# ---- Jason Schema sheet (A:Employee B:EmployeeID C:Date D:Client E:Hours F:Rate G:Total H:Type I:Notes) ----

$newEmpId = "emp_emnnysju"

# Row 2 - 2026-01-26 -> PTO day
$wsSchema.Range("B2").Value = $newEmpId
$wsSchema.Range("D2").Value = "PTO"
$wsSchema.Range("E2").Value = 6
$wsSchema.Range("F2").Value = 88
$wsSchema.Range("G2").Value = 528
$wsSchema.Range("H2").Value = "PTO"
$wsSchema.Range("I2").Value = "PTO"

# Row 3 - 2026-01-27 (Notes column I3 stays an untouched blank cell)
$wsSchema.Range("B3").Value = $newEmpId
$wsSchema.Range("D3").Value = "Vincent"
$wsSchema.Range("E3").Value = 6.5
$wsSchema.Range("F3").Value = 88
$wsSchema.Range("G3").Value = 572
$wsSchema.Range("H3").Value = "Regular"

# Row 4 - 2026-01-28 (Notes column I4 stays an untouched blank cell)
$wsSchema.Range("B4").Value = $newEmpId
$wsSchema.Range("D4").Value = "Zygmunt"
$wsSchema.Range("E4").Value = 6
$wsSchema.Range("F4").Value = 88
$wsSchema.Range("G4").Value = 528
$wsSchema.Range("H4").Value = "Regular"

# Row 5 - 2026-01-29 (Notes column I5 stays an untouched blank cell)
$wsSchema.Range("B5").Value = $newEmpId
$wsSchema.Range("D5").Value = "Ricca"
$wsSchema.Range("E5").Value = 6.5
$wsSchema.Range("F5").Value = 88
$wsSchema.Range("G5").Value = 572
$wsSchema.Range("H5").Value = "Regular"

# Row 6 - 2026-01-30 (Notes column I6 stays an untouched blank cell)
$wsSchema.Range("B6").Value = $newEmpId
$wsSchema.Range("D6").Value = "Varricchio"
$wsSchema.Range("E6").Value = 7
$wsSchema.Range("F6").Value = 88
$wsSchema.Range("G6").Value = 616
$wsSchema.Range("H6").Value = "Regular"
